# Update the "Eszközök" table header row on the Munka2 sheet to use
# colon-terminated labels, and add a new "Végösszeg:" row underneath
# the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Munka2")

# Header row (row 3): keep the same meaning, just append a colon to
# each label. Set in D3, E3, C3 order so new shared-string entries are
# created in that order (matches the original authoring order).
$ws.Range("D3").Value = "Egység ár:"
$ws.Range("E3").Value = "Össz. Ár:"
$ws.Range("C3").Value = "Darab szám:"

# New summary row below the table.
$ws.Range("D19").Value = "Végösszeg:"

# Leave the selection on the newly edited cell, as in the saved file.
$ws.Range("E19").Select()
